# Refine metadata to be an additional tab.
#
# The workbook currently has a single "data" sheet. This adds a new
# "metadata" sheet (after "data") describing the panel query that produced
# the "data" sheet, and refreshes the "time_taken" timestamps on "data" to
# match the new query run.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # "data"

# --- Add the new "metadata" worksheet after the last existing sheet -------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet  = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$newSheet.Name = "metadata"

# --- Reuse the existing bold/centered/bordered header style from "data" ---
# (copy formats only, so no new style entries are introduced)
$ws1.Range("F1").Copy()
$newSheet.Range("G1").PasteSpecial(-4122)   # xlPasteFormats

$ws1.Range("B1:F1").Copy()
$newSheet.Range("B1:F1").PasteSpecial(-4122)

# Reuse the numeric-index style (style of A2 on "data") for A2 on "metadata"
$ws1.Range("A2").Copy()
$newSheet.Range("A2").PasteSpecial(-4122)

# --- Header row -------------------------------------------------------
$newSheet.Range("B1").Value = "data_name"
$newSheet.Range("C1").Value = "data_id"
$newSheet.Range("D1").Value = "data_version"
$newSheet.Range("E1").Value = "data_version_created"
$newSheet.Range("F1").Value = "panel_query_time"
$newSheet.Range("G1").Value = "panel_get_request"

# --- Data row -----------------------------------------------------------
$newSheet.Range("A2").Value = 0
$newSheet.Range("B2").Value = "Hydroa vacciniforme"
$newSheet.Range("C2").Value = 310

# "data_version" must stay textual ("1.2"), not be coerced to a number.
# Mark the cell as Text before typing the value, then drop the formatting
# again so the cell ends up back on the default (unstyled) format.
$d2 = $newSheet.Range("D2")
$d2.NumberFormat = "@"
$d2.Value = "1.2"
$d2.ClearFormats()

$newSheet.Range("E2").Value = "2017-11-05T02:37:20.406396Z"
$newSheet.Range("F2").Value = "2021-10-05 14:20:50.214831"
$newSheet.Range("G2").Value = "https://panelapp.genomicsengland.co.uk/api/v1/panels/310/?format=json"

# --- Refresh the time_taken timestamps on "data" to match the new run ---
$ws1.Range("F2").Value = "2021-10-05 14:20:50.218554"
$ws1.Range("F3").Value = "2021-10-05 14:20:50.218562"
$ws1.Range("F4").Value = "2021-10-05 14:20:50.218565"
